$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, [string]$val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "24.647.92"
Set-TextValue "E2" "  +3.46%  "

Set-TextValue "D3" "1.701.00"

Set-TextValue "E4" "  -0.01%  "

Set-TextValue "D5" "316.12"
Set-TextValue "E5" "  +2.42%  "

Set-TextValue "D7" "0.3946"
Set-TextValue "E7" "  +1.06%  "

Set-TextValue "D8" "0.4043"
Set-TextValue "E8" "  +2.30%  "

Set-TextValue "D9" "1.542"
Set-TextValue "E9" "  +8.69%  "

Set-TextValue "D10" "55.26"
Set-TextValue "E10" "  +13.51%  "

Set-TextValue "E11" "  -0.03%  "

Set-TextValue "D12" "0.08822"
Set-TextValue "E12" "  +2.45%  "

Set-TextValue "D13" "7.310"
Set-TextValue "E13" "  +8.51%  "

Set-TextValue "E14" "  +3.48%  "

Set-TextValue "E15" "  +1.86%  "

Set-TextValue "D17" "1.704.05"
Set-TextValue "E17" "  +2.46%  "

Set-TextValue "D18" "100.98"
Set-TextValue "E18" "  +1.08%  "

Set-TextValue "D19" "0.07065"
Set-TextValue "E19" "  +4.22%  "

Set-TextValue "E20" "  +3.93%  "

Set-TextValue "D21" "6.938"
Set-TextValue "E21" "  +4.66%  "

Set-TextValue "D22" "1.000"
Set-TextValue "E22" "  -0.03%  "

Set-TextValue "D23" "14.19"
Set-TextValue "E23" "  +3.18%  "

Set-TextValue "D24" "24.631.96"
Set-TextValue "E24" "  +3.50%  "

Set-TextValue "D25" "2.996"
Set-TextValue "E25" "  +10.48%  "

Set-TextValue "E26" "  +0.28%  "

Set-TextValue "E27" "  +3.55%  "

Set-TextValue "D28" "160.10"
Set-TextValue "E28" "  +1.84%  "

Set-TextValue "D29" "5.236"
Set-TextValue "E29" "  +0.94%  "

Set-TextValue "D30" "134.13"
Set-TextValue "E30" "  +3.64%  "

Set-TextValue "D31" "7.771"
Set-TextValue "E31" "  +19.44%  "

Set-TextValue "D32" "1.112"
Set-TextValue "E32" "  -1.27%  "

Set-TextValue "D33" "1.888.73"
Set-TextValue "E33" "  +2.29%  "

Set-TextValue "D34" "7.441"
Set-TextValue "E34" "  +14.79%  "

Set-TextValue "D35" "0.08588"
Set-TextValue "E35" "  -0.13%  "

Set-TextValue "D36" "11.21"
Set-TextValue "E36" "  +8.29%  "

Set-TextValue "D37" "0.2764"
Set-TextValue "E37" "  +4.91%  "

Set-TextValue "D38" "1.952"
Set-TextValue "E38" "  -0.13%  "

Set-TextValue "D39" "14.83"
Set-TextValue "E39" "  +2.83%  "

Set-TextValue "D40" "0.02797"
Set-TextValue "E40" "  +10.95%  "

Set-TextValue "D41" "0.09060"
Set-TextValue "E41" "  +3.45%  "

Set-TextValue "D42" "1.474"
Set-TextValue "E42" "  +2.73%  "

Set-TextValue "D43" "0.7773"
Set-TextValue "E43" "  +3.46%  "

Set-TextValue "D44" "0.7297"
Set-TextValue "E44" "  +4.26%  "

Set-TextValue "D45" "15.61"
Set-TextValue "E45" "  +5.11%  "

Set-TextValue "D46" "2.521"
Set-TextValue "E46" "  +6.70%  "

Set-TextValue "D47" "4.208"
Set-TextValue "E47" "  +3.43%  "

Set-TextValue "B48" "Flow"
Set-TextValue "C48" "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
Set-TextValue "D48" "1.391"
Set-TextValue "E48" "  +21.28%  "

Set-TextValue "B49" "Frax"
Set-TextValue "C49" "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D49" "0.9999"
Set-TextValue "E49" "  -0.07%  "

Set-TextValue "B50" "Quant"
Set-TextValue "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "141.91"
Set-TextValue "E50" "  +1.78%  "

Set-TextValue "D51" "0.08037"
Set-TextValue "E51" "  +3.72%  "
